$d = $word.ActiveDocument

# ---- 1. Remove _GoBack bookmark from start of doc ----
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---- 2. Version bump 5.10.0 -> 5.10.1 in title ----
$d.Content.Find.Execute(".0 Release Notes", $true, $false, $false, $false, $false, $true, 1, $false, ".1 Release Notes", 2) | Out-Null

# ---- 3. Date change ----
$d.Content.Find.Execute("January 13, 2023", $true, $false, $false, $false, $false, $true, 1, $false, "March 20, 2023", 2) | Out-Null

# ---- 4. Current Release Notes (5.10.0) -> (5.10.1) ----
$d.Content.Find.Execute("Current Release Notes (5.10.0)", $true, $false, $false, $false, $false, $true, 1, $false, "Current Release Notes (5.10.1)", 2) | Out-Null

# ---- 5. Insert bold " (5.10.0)" after "table has been converted to SQLite" ----
$rng = $d.Content
$rng.Find.Execute("table has been converted to SQLite") | Out-Null
$insPt = $d.Range($rng.End, $rng.End)
$insPt.InsertAfter(" (5.10.0)")
$insPt.Font.Bold = 1

# ---- 6. Move _GoBack bookmark to after "Please contact BioSum support for instructions." ----
$rng2 = $d.Content
$rng2.Find.Execute("Please contact BioSum support for instructions.") | Out-Null
$bmPt = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $bmPt)

# ---- 7. Insert bold " (5.10.0)" after "Updates to FVS sequence number definition screen" ----
$rng3 = $d.Content
$rng3.Find.Execute("Updates to FVS sequence number definition screen") | Out-Null
$insPt3 = $d.Range($rng3.End, $rng3.End)
$insPt3.InsertAfter(" (5.10.0)")
$insPt3.Font.Bold = 1

# ---- 8. Insert a new list paragraph before "Many additional minor enhancements..." ----
# (paragraph 5 is the "Many additional minor enhancements..." paragraph at this point)
$p5 = $d.Paragraphs.Item(5)
$insBefore = $d.Range($p5.Range.Start, $p5.Range.Start)
$insBefore.InsertParagraphBefore()

# The freshly-inserted paragraph is now paragraph 5 (it inherits the ListParagraph /
# numbering format from the paragraph it was inserted in front of).
$newPara = $d.Paragraphs.Item(5)
$paraStart = $newPara.Range.Start

$boldText = "Additional harvest costs can now be assessed using variable in the KCP file (5.10.1)"
$restText = ": Additional harvest costs can be assessed at the Rx/RxYear level by configuring flags in the FVS_Compute table using KCP file directives. The dollar amount of the costs is configured in the Processor module. Please see the user guide for details."

$ip = $d.Range($paraStart, $paraStart)
$ip.InsertAfter($boldText + $restText)

$boldRange = $d.Range($paraStart, $paraStart + $boldText.Length)
$boldRange.Font.Bold = 1
$restRange = $d.Range($paraStart + $boldText.Length, $paraStart + $boldText.Length + $restText.Length)
$restRange.Font.Bold = 0

$brPart = $d.Range($restRange.End, $restRange.End)
$brPart.InsertAfter([char]11)

Write-Host "done"
